# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 320
$ws1.Range("F4").Value = 1292
$ws1.Range("F5").Value = 638

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 320
$ws4.Range("F4").Value = 1292
$ws4.Range("F6").Value = 638
